# The edit reshuffles the per-observation rows (2-14) on the "Artfynd"
# sheet: each row's whole record (Id, taxon info, coordinates, comments,
# ...) moves to a different row, per the mapping derived from the diff
# (new row -> source row holding its old content):
#   2<-10  3<-5  4<-2  5<-4  6<-9  7<-12  8<-8  9<-14  10<-7  11<-13  12<-6  13<-11  14<-3
# Row 8 maps to itself (unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 14
$rangeAddr = "A" + $firstRow + ":AY" + $lastRow

# Startdatum/Slutdatum (Y, AA) hold plain "yyyy-mm-dd" text in the source
# file, not real dates. Force text format up front so the Value2
# round-trip below doesn't silently promote them to date serials.
$ws.Range("Y" + $firstRow + ":Y" + $lastRow).NumberFormat = "@"
$ws.Range("AA" + $firstRow + ":AA" + $lastRow).NumberFormat = "@"

# Snapshot every value in the data block before writing anything back.
$snapshot = $ws.Range($rangeAddr).Value2

# new row (key) <- source row (value), both are worksheet row numbers
$mapping = @{
    2  = 10
    3  = 5
    4  = 2
    5  = 4
    6  = 9
    7  = 12
    8  = 8
    9  = 14
    10 = 7
    11 = 13
    12 = 6
    13 = 11
    14 = 3
}

# $snapshot (read off a Range) is 1-based: [1..rowCount, 1..colCount].
# A freshly allocated .NET array is 0-based, and that is what Excel
# expects back on assignment - so source/dest indices differ by one.
$rowCount = $lastRow - $firstRow + 1
$colCount = $snapshot.GetLength(1)

$result = New-Object 'object[,]' $rowCount, $colCount

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapping[$r]
    $destIdx = $r - $firstRow       # 0-based
    $srcIdx = $srcRow - $firstRow + 1  # 1-based (into $snapshot)
    for ($c = 1; $c -le $colCount; $c++) {
        $result[$destIdx, $c - 1] = $snapshot[$srcIdx, $c]
    }
}

$ws.Range($rangeAddr).Value2 = $result

# --- Reproduce the handful of cells that are *present but empty* (an
# inline string with no text) in the target layout, which the bulk
# Value2 copy above cannot express (it only carries real values, so an
# empty-but-present cell and an absent cell both come back blank).
# Assigning a NumberFormat to an empty cell makes Excel materialise it
# in the sheet XML without giving it a value.
#
# Columns I, AT and AY are blank-but-present on every data row.
# Column AF is blank-but-present only on the rows whose (permuted)
# content originally carried that marker: new rows 2, 3, 7, 11.
$blankMarkerCells = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $blankMarkerCells += "I" + $r
    $blankMarkerCells += "AT" + $r
    $blankMarkerCells += "AY" + $r
}
$blankMarkerCells += @("AF2", "AF3", "AF7", "AF11")

foreach ($addr in $blankMarkerCells) {
    $ws.Range($addr).NumberFormat = "@"
}
